$d = $word.ActiveDocument

# Locate the paragraph that holds the field (fldChar begin/instrText/.../fldChar end).
# The field's "m:'doc.html'.fromHTMLURI()" M2Doc query is currently stored as a real
# Word field; the commit turns it into plain literal text "{ m:'doc.html'.fromHTMLURI() }"
# (braces replacing the field delimiters) so the new TokenIteratorFieldRewriterSplit
# parser can tokenize it directly from the run text, while keeping the _GoBack
# bookmark exactly where it was.
$fieldPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $fieldPara = $p
    }
}

if ($fieldPara -ne $null) {
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">' +
           '<w:r><w:t>{</w:t></w:r>' +
           '<w:r><w:t>m</w:t></w:r>' +
           '<w:r><w:t>:</w:t></w:r>' +
           "<w:r><w:t>'</w:t></w:r>" +
           '<w:r><w:t>doc.html</w:t></w:r>' +
           '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
           '<w:bookmarkEnd w:id="0"/>' +
           "<w:r><w:t>'.fromHTMLURI()</w:t></w:r>" +
           '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
           '</w:p>'
    [void]$fieldPara.Range.InsertXML($xml)
}
